$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Block at row 66: previously a placeholder header "Dec 16" with no data.
# Relabel it to "Dec 22" and fill in the measured input/output values.
$ws.Range("A66").Value = "Dec 22"

$ws.Range("B67").Value = 77.6
$ws.Range("C67").Value = 89.7
$ws.Range("D67").Value = 85
$ws.Range("E67").Value = 86.1

$ws.Range("B68").Value = 78
$ws.Range("C68").Value = 86.6
$ws.Range("D68").Value = 90.4
$ws.Range("E68").Value = 88.9

$ws.Range("B69").Value = 75.6
$ws.Range("C69").Value = 85.7
$ws.Range("D69").Value = 83.2
$ws.Range("E69").Value = 76.5

$ws.Range("B70").Value = 73.1
$ws.Range("C70").Value = 89.3
$ws.Range("D70").Value = 82.7
$ws.Range("E70").Value = 81.5

$ws.Range("B71").Value = 83.8
$ws.Range("C71").Value = 88.1
$ws.Range("D71").Value = 89.7
$ws.Range("E71").Value = 88.3

# --- Block at row 74: previously a placeholder header "Dec 20" with no data.
# Relabel it to "Dec 23" and fill in the measured input/output values.
$ws.Range("A74").Value = "Dec 23"

$ws.Range("B75").Value = 79.2
$ws.Range("C75").Value = 89.8
$ws.Range("D75").Value = 86.9
$ws.Range("E75").Value = 86.9

$ws.Range("B76").Value = 79.6
$ws.Range("C76").Value = 86.2
$ws.Range("D76").Value = 91
$ws.Range("E76").Value = 90.3

$ws.Range("B77").Value = 77
$ws.Range("C77").Value = 84.2
$ws.Range("D77").Value = 83.9
$ws.Range("E77").Value = 77.8

$ws.Range("B78").Value = 73.7
$ws.Range("C78").Value = 89.9
$ws.Range("D78").Value = 83.6
$ws.Range("E78").Value = 80.5

$ws.Range("B79").Value = 85
$ws.Range("C79").Value = 92.6
$ws.Range("D79").Value = 90.1
$ws.Range("E79").Value = 87.6

# --- Block at rows 82-87: this was a leftover empty placeholder block
# ("Dec 23" header, no data) that's no longer needed now that the data
# above has been filled in. Delete those rows entirely.
$ws.Rows("82:87").Delete()

# Update the visible scroll/selection position to match where the user
# was working after the edit.
$ws.Application.ActiveWindow.ScrollRow = 47
$ws.Range("E80").Select()
